$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 0.9999919691916495
$ws.Range("E2").Value = 0.9999919691916495

# Row 3
$ws.Range("D3").Value = 0.9999998948690316
$ws.Range("E3").Value = 0.9999998948690316

# Row 4
$ws.Range("D4").Value = [double]"1.149104031854136E-21"
$ws.Range("E4").Value = [double]"1.149104031854136E-21"

# Row 5
$ws.Range("D5").Value = 0.0001429388540569712
$ws.Range("E5").Value = 0.0001429388540569712

# Row 6
$ws.Range("C6").Value = $true
$ws.Range("D6").Value = 0.2450349435315865
$ws.Range("E6").Value = 0.2450349435315865

# Row 7
$ws.Range("D7").Value = 0.999999999796088
$ws.Range("E7").Value = [double]"2.039119983976434E-10"

# Row 8
$ws.Range("D8").Value = 0.9999999999966604
$ws.Range("E8").Value = [double]"3.339550858072471E-12"

# Row 9
$ws.Range("D9").Value = 0.0006436159059134271
$ws.Range("E9").Value = 0.9993563840940866

# Row 10
$ws.Range("D10").Value = 0.9999999999999782
$ws.Range("E10").Value = [double]"2.176037128265307E-14"

# Row 11
$ws.Range("D11").Value = [double]"8.292513112645639E-06"
$ws.Range("E11").Value = 0.9999917074868874
$ws.Range("F11").Value = 4.71300745010376
$ws.Range("G11").Value = 0.6
